$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the systematic/statistical uncertainty values (columns G, H, I) and
# a couple of corrected "diff" values (D7, D8) for rows 2-10.
$ws.Range("D7").Value = 0.37
$ws.Range("D8").Value = 0.64

$ws.Range("G2").Value = 1.9
$ws.Range("H2").Value = 0.2
$ws.Range("I2").Value = 0.1

$ws.Range("G3").Value = 3.8
$ws.Range("H3").Value = 0.3
$ws.Range("I3").Value = 0.1

$ws.Range("G4").Value = 4.6
$ws.Range("H4").Value = 0.3
$ws.Range("I4").Value = 0.1

$ws.Range("G5").Value = 5.9
$ws.Range("H5").Value = 0.4
$ws.Range("I5").Value = 0.2

$ws.Range("G6").Value = 6.7
$ws.Range("H6").Value = 0.5
$ws.Range("I6").Value = 0.3

$ws.Range("G7").Value = 5.4
$ws.Range("H7").Value = 0.4
$ws.Range("I7").Value = 0.3

$ws.Range("G8").Value = 3.3
$ws.Range("I8").Value = 0.1

$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 0.2
$ws.Range("I9").Value = 0.1

$ws.Range("G10").Value = 0.9
$ws.Range("H10").Value = 0.1
$ws.Range("I10").Value = 0.1

# Remove the stray leftover cell/row below the data table.
$ws.Rows("19:19").Delete()

# Restore the selection to match the saved workbook state.
$ws.Range("J17").Select()
